# Scheduled runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# figures (columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to
# reflect newly pulled market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 265
$ws.Range("I2").Value = 265
$ws.Range("K2").Value = 265
$ws.Range("M2").Value = -152
$ws.Range("H38").Value = 717237.9
$ws.Range("I38").Value = 1111392.4
$ws.Range("J38").Value = 7759.8
$ws.Range("K38").Value = 3334177.2
$ws.Range("L38").Value = 23279.4
$ws.Range("M38").Value = -3333805.2
$ws.Range("N38").Value = -24023.4
$ws.Range("H51").Value = 5660.467
$ws.Range("I51").Value = 5527.25
$ws.Range("K51").Value = 5527.25
$ws.Range("M51").Value = -5043.25
$ws.Range("H92").Value = 721.8
$ws.Range("I92").Value = 614.875
$ws.Range("K92").Value = 614.875
$ws.Range("M92").Value = 633.125
$ws.Range("H129").Value = 2361.1
$ws.Range("I129").Value = 2132
$ws.Range("J129").Value = 2386.5557
$ws.Range("K129").Value = 6396
$ws.Range("L129").Value = 7159.6671
$ws.Range("M129").Value = -1396
$ws.Range("N129").Value = -17159.6671
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 432.41666
$ws.Range("I5").Value = 511.42856
$ws.Range("J5").Value = 321.8
$ws.Range("K5").Value = 511.42856
$ws.Range("L5").Value = 321.8
$ws.Range("M5").Value = -399.42856
$ws.Range("N5").Value = -545.8
$ws.Range("H80").Value = 39988.375
$ws.Range("J80").Value = 39988.375
$ws.Range("L80").Value = 39988.375
$ws.Range("N80").Value = -41984.375
$ws.Range("H83").Value = 39988.375
$ws.Range("J83").Value = 39988.375
$ws.Range("L83").Value = 119965.125
$ws.Range("N83").Value = -129949.125
$ws.Range("H95").Value = 56949.332
$ws.Range("J95").Value = 56949.332
$ws.Range("L95").Value = 56949.332
$ws.Range("N95").Value = -62441.332
$ws.Range("H97").Value = 463.85715
$ws.Range("I97").Value = 482.83334
$ws.Range("K97").Value = 482.83334
$ws.Range("M97").Value = 13.16665999999998
$ws.Range("H102").Value = 2057.25
$ws.Range("I102").Value = 2057.25
$ws.Range("K102").Value = 2057.25
$ws.Range("M102").Value = -435.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 432.41666
$ws.Range("I4").Value = 511.42856
$ws.Range("J4").Value = 321.8
$ws.Range("K4").Value = 511.42856
$ws.Range("L4").Value = 321.8
$ws.Range("M4").Value = -396.42856
$ws.Range("N4").Value = -551.8
$ws.Range("H86").Value = 3199.5
$ws.Range("I86").Value = 5750
$ws.Range("J86").Value = 1924.25
$ws.Range("K86").Value = 5750
$ws.Range("L86").Value = 1924.25
$ws.Range("M86").Value = -4627
$ws.Range("N86").Value = -4170.25
$ws.Range("H89").Value = 3199.5
$ws.Range("I89").Value = 5750
$ws.Range("J89").Value = 1924.25
$ws.Range("K89").Value = 28750
$ws.Range("L89").Value = 9621.25
$ws.Range("M89").Value = -23134
$ws.Range("N89").Value = -20853.25
$ws.Range("H94").Value = 337.22223
$ws.Range("I94").Value = 337.22223
$ws.Range("K94").Value = 337.22223
$ws.Range("M94").Value = 113.77777
$ws.Range("H134").Value = 2516.625
$ws.Range("I134").Value = 2516.625
$ws.Range("K134").Value = 7549.875
$ws.Range("M134").Value = -5014.875
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 166667600
$ws.Range("I16").Value = 166667600
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 166667600
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 3811.8235
$ws.Range("I31").Value = 3070.5
$ws.Range("K31").Value = 3070.5
$ws.Range("M31").Value = -2775.5
$ws.Range("H34").Value = 3811.8235
$ws.Range("I34").Value = 3070.5
$ws.Range("K34").Value = 3070.5
$ws.Range("M34").Value = -2868.5
$ws.Range("H74").Value = 25385.666
$ws.Range("J74").Value = 24062.8
$ws.Range("L74").Value = 24062.8
$ws.Range("N74").Value = -25810.8
$ws.Range("H77").Value = 25385.666
$ws.Range("J77").Value = 24062.8
$ws.Range("L77").Value = 72188.39999999999
$ws.Range("N77").Value = -80924.39999999999
$ws.Range("H105").Value = 1228.5714
$ws.Range("J105").Value = 1283.3334
$ws.Range("L105").Value = 1283.3334
$ws.Range("N105").Value = -4777.3334
$ws.Range("H107").Value = 1294.3914
$ws.Range("I107").Value = 1192.5294
$ws.Range("K107").Value = 1192.5294
$ws.Range("M107").Value = 727.4706000000001
$ws.Range("H113").Value = 166667600
$ws.Range("I113").Value = 166667600
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 166667600
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 910
$ws.Range("I29").Value = 937.5
$ws.Range("J29").Value = 800
$ws.Range("K29").Value = 2812.5
$ws.Range("L29").Value = 2400
$ws.Range("M29").Value = -2535.5
$ws.Range("N29").Value = -2954
$ws.Range("H33").Value = 666.55554
$ws.Range("I33").Value = 671.2857
$ws.Range("J33").Value = 650
$ws.Range("K33").Value = 4027.7142
$ws.Range("L33").Value = 3900
$ws.Range("M33").Value = -3744.7142
$ws.Range("N33").Value = -4466
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("H131").Value = 2449.75
$ws.Range("I131").Value = 2359.8
$ws.Range("K131").Value = 7079.400000000001
$ws.Range("M131").Value = -2039.400000000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 396.69232
$ws.Range("I97").Value = 350.63635
$ws.Range("J97").Value = 650
$ws.Range("K97").Value = 350.63635
$ws.Range("L97").Value = 650
$ws.Range("M97").Value = 145.36365
$ws.Range("N97").Value = -1642
$ws.Range("H113").Value = 1300.8334
$ws.Range("I113").Value = 1191.8182
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1191.8182
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 978.1818000000001
$ws.Range("N113").Value = -6840
$ws.Range("H122").Value = 2336
$ws.Range("I122").Value = 2179.6
$ws.Range("K122").Value = 6538.799999999999
$ws.Range("M122").Value = -4088.799999999999
$ws.Range("H126").Value = 1798
$ws.Range("J126").Value = 1798
$ws.Range("L126").Value = 5394
$ws.Range("N126").Value = -10334
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4426.6924
$ws.Range("I82").Value = 2949.4285
$ws.Range("J82").Value = 6150.1665
$ws.Range("K82").Value = 2949.4285
$ws.Range("L82").Value = 6150.1665
$ws.Range("M82").Value = -2588.4285
$ws.Range("N82").Value = -6872.1665
$ws.Range("H85").Value = 4426.6924
$ws.Range("I85").Value = 2949.4285
$ws.Range("J85").Value = 6150.1665
$ws.Range("K85").Value = 2949.4285
$ws.Range("L85").Value = 6150.1665
$ws.Range("M85").Value = -1701.4285
$ws.Range("N85").Value = -8646.166499999999
$ws.Range("H122").Value = 16999.875
$ws.Range("J122").Value = 16333
$ws.Range("L122").Value = 48999
$ws.Range("N122").Value = -53899
$ws.Range("H132").Value = 2968.7222
$ws.Range("I132").Value = 2962.4666
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 8887.399800000001
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -6357.399800000001
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 45457396
$ws.Range("I136").Value = 4058.1667
$ws.Range("K136").Value = 12174.5001
$ws.Range("M136").Value = -9624.500100000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5899
$ws.Range("J96").Value = 5899
$ws.Range("L96").Value = 5899
$ws.Range("N96").Value = -8645
$ws.Range("H122").Value = 2207.5
$ws.Range("I122").Value = 1677.8948
$ws.Range("J122").Value = 4220
$ws.Range("K122").Value = 5033.6844
$ws.Range("L122").Value = 12660
$ws.Range("M122").Value = -2583.6844
$ws.Range("N122").Value = -17560
$ws.Range("H132").Value = 5264.9165
$ws.Range("I132").Value = 4662.5
$ws.Range("K132").Value = 13987.5
$ws.Range("M132").Value = -11457.5
